$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename strategy labels in column B from E1..E7 to S1..S7
for ($i = 1; $i -le 7; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = "S$i"
}
